$wb = $excel.ActiveWorkbook

# "About" sheet: remove the stray date value in C1 (leftover stamp)
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Clear()

# "BENCEfCT" sheet: flip the boolean lever value in B2 from 1 to 0
$wsLever = $wb.Worksheets.Item("BENCEfCT")
$wsLever.Range("B2").Value = 0

# Make "About" the active/selected sheet (was BENCEfCT)
$wsAbout.Activate()
